# Append a new row (row 26) of bitcoin purchase data, matching the
# automated run performed on 2025-06-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 26

# Column A: the date is stored as literal text "06/25/2025" (same
# convention as the other recently-appended rows), not a date serial.
# Force text formatting so Excel does not auto-parse the string into a
# date value, then restore the default "Normal" style so the cell ends
# up with no explicit style override (matching the existing rows).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "06/25/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

# Column B: coins purchased
$ws.Cells.Item($row, 2).Value = 0.0004631500000000007

# Column C: price
$ws.Cells.Item($row, 3).Value = 107956.3856202093

# Column D: cost
$ws.Cells.Item($row, 4).Value = 50
